# Acabamos el tema de Mybatis con Manolo Empezamos con Hibernate 17/04/2018 11:45
#
# Appends 6 new document rows (doc name "nombre", status "ACTIVO") with a
# series of timestamps to sheet "Hoja3" (rows 12-17) and to sheet "Hoja2"
# (rows 6-11).

$wb = $excel.ActiveWorkbook

$timestamps = @(
    "Tue Apr 17 09:52:40 CEST 2018",
    "Tue Apr 17 09:52:55 CEST 2018",
    "Tue Apr 17 09:53:27 CEST 2018",
    "Tue Apr 17 09:59:00 CEST 2018",
    "Tue Apr 17 10:01:15 CEST 2018",
    "Tue Apr 17 10:01:22 CEST 2018"
)

# --- Hoja3: currently rows 2-11 populated (A2:F11) -> add rows 12-17 ---
$wsHoja3 = $wb.Worksheets.Item("Hoja3")
$startRow = 12
for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $startRow + $i
    $wsHoja3.Cells.Item($row, 1).Value = $row - 1
    $wsHoja3.Cells.Item($row, 2).Value = 1
    $wsHoja3.Cells.Item($row, 3).Value = "nombre"
    $wsHoja3.Cells.Item($row, 4).Value = $timestamps[$i]
    $wsHoja3.Cells.Item($row, 5).Value = "ACTIVO"
}

# --- Hoja2: currently rows 2-5 populated (A2:F5) -> add rows 6-11 ---
$wsHoja2 = $wb.Worksheets.Item("Hoja2")
$startRow = 6
for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $startRow + $i
    $wsHoja2.Cells.Item($row, 1).Value = $row - 1
    $wsHoja2.Cells.Item($row, 2).Value = 1
    $wsHoja2.Cells.Item($row, 3).Value = "nombre"
    $wsHoja2.Cells.Item($row, 4).Value = $timestamps[$i]
    $wsHoja2.Cells.Item($row, 5).Value = "ACTIVO"
}
